$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# Update Testcase ID values in column A for rows 10-13
$ws.Range("A10").Value = 9
$ws.Range("A11").Value = 10
$ws.Range("A12").Value = 11
$ws.Range("A13").Value = 12

# Update the view: scroll so row 12 is the top-left visible row,
# and select H12
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Range("H12").Select()
